$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 126, pushing the existing rows 126-143
# down to 128-145 (dimension grows from A1:T143 to A1:T145).
$ws.Range("A126:A127").EntireRow.Insert()

# --- New row 126 : Provincia de Curicó, Primera ---
$ws.Cells.Item(126, 1).Value = 9
$ws.Cells.Item(126, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(126, 3).Value = "Metropolitana"
$ws.Cells.Item(126, 4).Value = 44984
$ws.Cells.Item(126, 5).Value = 13
$ws.Cells.Item(126, 6).Value = "Fruta"
$ws.Cells.Item(126, 7).Value = 100101
$ws.Cells.Item(126, 8).Value = "Berries"
$ws.Cells.Item(126, 9).Value = 100101004
$ws.Cells.Item(126, 10).Value = "Frambuesa"
$ws.Cells.Item(126, 11).Value = "Sin especificar"
$ws.Cells.Item(126, 12).Value = "Primera"
$ws.Cells.Item(126, 13).Value = 180
$ws.Cells.Item(126, 14).Value = 6000
$ws.Cells.Item(126, 15).Value = 6000
$ws.Cells.Item(126, 16).Value = 6000
$ws.Cells.Item(126, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(126, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(126, 19).Value = 3000
$ws.Cells.Item(126, 20).Value = 2

# --- New row 127 : Provincia de Curicó, Segunda ---
$ws.Cells.Item(127, 1).Value = 9
$ws.Cells.Item(127, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(127, 3).Value = "Metropolitana"
$ws.Cells.Item(127, 4).Value = 44984
$ws.Cells.Item(127, 5).Value = 13
$ws.Cells.Item(127, 6).Value = "Fruta"
$ws.Cells.Item(127, 7).Value = 100101
$ws.Cells.Item(127, 8).Value = "Berries"
$ws.Cells.Item(127, 9).Value = 100101004
$ws.Cells.Item(127, 10).Value = "Frambuesa"
$ws.Cells.Item(127, 11).Value = "Sin especificar"
$ws.Cells.Item(127, 12).Value = "Segunda"
$ws.Cells.Item(127, 13).Value = 150
$ws.Cells.Item(127, 14).Value = 5000
$ws.Cells.Item(127, 15).Value = 5000
$ws.Cells.Item(127, 16).Value = 5000
$ws.Cells.Item(127, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(127, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(127, 19).Value = 2500
$ws.Cells.Item(127, 20).Value = 2

# Make sure the date cells keep the same display format used by the
# other rows in column D (style index 2 -> "YYYY-MM-DD HH:MM:SS").
$ws.Cells.Item(126, 4).NumberFormat = $ws.Cells.Item(128, 4).NumberFormat
$ws.Cells.Item(127, 4).NumberFormat = $ws.Cells.Item(128, 4).NumberFormat
